{"js": "// Update the date line and the 25 division-problem answers in the table.\n// Each \"from\" text is unique within the document, so a simple\n// search-and-replace (by exact, case-sensitive match) for each pair is safe.\nconst replacements = [\n  [\"2025-08-22 Friday\", \"2025-08-23 Saturday\"],\n  [\"406\u00f75=81, 1\", \"569\u00f74=142, 1\"],\n  [\"810\u00f73=270, 0\", \"965\u00f77=137, 6\"],\n  [\"549\u00f78=68, 5\", \"449\u00f77=64, 1\"],\n  [\"212\u00f76=35, 2\", \"434\u00f74=108, 2\"],\n  [\"917\u00f79=101, 8\", \"268\u00f78=33, 4\"],\n  [\"273\u00f74=68, 1\", \"285\u00f74=71, 1\"],\n  [\"672\u00f72=336, 0\", \"999\u00f76=166, 3\"],\n  [\"707\u00f78=88, 3\", \"435\u00f79=48, 3\"],\n  [\"363\u00f78=45, 3\", \"981\u00f79=109, 0\"],\n  [\"311\u00f77=44, 3\", \"907\u00f72=453, 1\"],\n  [\"625\u00f74=156, 1\", \"520\u00f73=173, 1\"],\n  [\"991\u00f77=141, 4\", \"105\u00f76=17, 3\"],\n  [\"295\u00f78=36, 7\", \"990\u00f78=123, 6\"],\n  [\"858\u00f74=214, 2\", \"172\u00f79=19, 1\"],\n  [\"633\u00f75=126, 3\", \"955\u00f79=106, 1\"],\n  [\"803\u00f75=160, 3\", \"791\u00f75=158, 1\"],\n  [\"144\u00f74=36, 0\", \"551\u00f76=91, 5\"],\n  [\"479\u00f78=59, 7\", \"853\u00f73=284, 1\"],\n  [\"480\u00f72=240, 0\", \"702\u00f75=140, 2\"],\n  [\"607\u00f79=67, 4\", \"552\u00f73=184, 0\"],\n  [\"321\u00f74=80, 1\", \"259\u00f73=86, 1\"],\n  [\"168\u00f75=33, 3\", \"697\u00f72=348, 1\"],\n  [\"202\u00f78=25, 2\", \"195\u00f77=27, 6\"],\n  [\"108\u00f75=21, 3\", \"861\u00f79=95, 6\"],\n  [\"966\u00f74=241, 2\", \"128\u00f79=14, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 division-problem answers in the table.\n# Each \"Find\" text is unique within the document, so a straightforward\n# Find/Replace (exact match, case-sensitive) for each pair is safe.\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @{ Find = \"2025-08-22 Friday\"; Replace = \"2025-08-23 Saturday\" },\n    @{ Find = \"406\u00f75=81, 1\"; Replace = \"569\u00f74=142, 1\" },\n    @{ Find = \"810\u00f73=270, 0\"; Replace = \"965\u00f77=137, 6\" },\n    @{ Find = \"549\u00f78=68, 5\"; Replace = \"449\u00f77=64, 1\" },\n    @{ Find = \"212\u00f76=35, 2\"; Replace = \"434\u00f74=108, 2\" },\n    @{ Find = \"917\u00f79=101, 8\"; Replace = \"268\u00f78=33, 4\" },\n    @{ Find = \"273\u00f74=68, 1\"; Replace = \"285\u00f74=71, 1\" },\n    @{ Find = \"672\u00f72=336, 0\"; Replace = \"999\u00f76=166, 3\" },\n    @{ Find = \"707\u00f78=88, 3\"; Replace = \"435\u00f79=48, 3\" },\n    @{ Find = \"363\u00f78=45, 3\"; Replace = \"981\u00f79=109, 0\" },\n    @{ Find = \"311\u00f77=44, 3\"; Replace = \"907\u00f72=453, 1\" },\n    @{ Find = \"625\u00f74=156, 1\"; Replace = \"520\u00f73=173, 1\" },\n    @{ Find = \"991\u00f77=141, 4\"; Replace = \"105\u00f76=17, 3\" },\n    @{ Find = \"295\u00f78=36, 7\"; Replace = \"990\u00f78=123, 6\" },\n    @{ Find = \"858\u00f74=214, 2\"; Replace = \"172\u00f79=19, 1\" },\n    @{ Find = \"633\u00f75=126, 3\"; Replace = \"955\u00f79=106, 1\" },\n    @{ Find = \"803\u00f75=160, 3\"; Replace = \"791\u00f75=158, 1\" },\n    @{ Find = \"144\u00f74=36, 0\"; Replace = \"551\u00f76=91, 5\" },\n    @{ Find = \"479\u00f78=59, 7\"; Replace = \"853\u00f73=284, 1\" },\n    @{ Find = \"480\u00f72=240, 0\"; Replace = \"702\u00f75=140, 2\" },\n    @{ Find = \"607\u00f79=67, 4\"; Replace = \"552\u00f73=184, 0\" },\n    @{ Find = \"321\u00f74=80, 1\"; Replace = \"259\u00f73=86, 1\" },\n    @{ Find = \"168\u00f75=33, 3\"; Replace = \"697\u00f72=348, 1\" },\n    @{ Find = \"202\u00f78=25, 2\"; Replace = \"195\u00f77=27, 6\" },\n    @{ Find = \"108\u00f75=21, 3\"; Replace = \"861\u00f79=95, 6\" },\n    @{ Find = \"966\u00f74=241, 2\"; Replace = \"128\u00f79=14, 2\" }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($item in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $item.Find\n    $find.Replacement.Text = $item.Replace\n    $find.Forward = $true\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $find.MatchSoundsLike, $find.MatchAllWordForms, $find.Forward, $wdFindContinue, $find.Format, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n}\n"}
